# Update the "K" column (column G) values in the save_data sheet.
# The source stats were regenerated (K instead of Strike#, regen std/mean,
# calc and write s_vals), which changed the strikeout counts recorded for
# each outing. Rows 28 and 29 are unaffected (values stay 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 1
    9  = 2
    10 = 3
    11 = 1
    12 = 2
    13 = 1
    14 = 4
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    30 = 1
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
